$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replacement Hiring % (column G) values each decreased by 0.5 for rows 2-23
$ws.Range("G2").Value = 1.5
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 3.0000000000000004
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2.5
$ws.Range("G7").Value = 3.0000000000000004
$ws.Range("G8").Value = 2.5
$ws.Range("G9").Value = 3.5
$ws.Range("G10").Value = 3.5
$ws.Range("G11").Value = 5
$ws.Range("G12").Value = 3.0000000000000004
$ws.Range("G13").Value = 2
$ws.Range("G14").Value = 4
$ws.Range("G15").Value = 2.5
$ws.Range("G16").Value = 1.5
$ws.Range("G17").Value = 3.0000000000000004
$ws.Range("G18").Value = 4
$ws.Range("G19").Value = 2.5
$ws.Range("G20").Value = 2.5
$ws.Range("G21").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("G23").Value = 3.0000000000000004

# Update the saved selection state to match the authored view (F28)
$ws.Range("F28").Select()
